$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.763.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "'2.161.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.90%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'227.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'0.627"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'63.31"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "'0.0845"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "'15.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "'2.483.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").Value = "'21.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "'0.806"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "'2.148.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "'39.639.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").Value = "'71.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "'0.0₃0845"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "'229.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").Value = "'2.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.07%  "
$ws.Range("D26").Value = "'172.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "'9.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("D30").Value = "'19.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").Value = "'2.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.21%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("D34").Value = "'4.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("D35").Value = "'6.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("D36").Value = "'0.0617"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("E37").Value = "  +4.81%  "
$ws.Range("E38").Value = "  +1.85%  "
$ws.Range("D39").Value = "'5.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +22.37%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "'102.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("D42").Value = "'17.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").Value = "'1.515.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "'0.0919"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("E50").Value = "  +8.73%  "
$ws.Range("D51").Value = "'3.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.28%  "
